$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns L (12) through S (19)
$headers = @(
    "hzj-混合调节_20170516_152754_ASIC_EEG",
    "hzj-混合调节_20170518_134207_ASIC_EEG",
    "hzj-混合调节_20170519_135415_ASIC_EEG",
    "zyx-混合调节_20170516_111228_ASIC_EEG",
    "zyx-混合调节_20170517_110944_ASIC_EEG",
    "zyx-混合调节_20170518_112337_ASIC_EEG",
    "zyx-混合调节_20170519_124954_ASIC_EEG",
    "zyx-混合调节_20170522_111557_ASIC_EEG"
)

# Row 2 values for columns L through S
$row2 = @(
    0.71875,
    0.74914089347079038,
    0.76029962546816476,
    0.74038461538461542,
    0.83333333333333337,
    0.78640776699029125,
    0.73225806451612896,
    0.78883495145631066
)

# Row 3 values for columns L through S
$row3 = @(
    0.69230769230769229,
    0.70322580645161292,
    0.69597069597069594,
    0.70714285714285707,
    0.74301675977653625,
    0.72852233676975953,
    0.74496644295302006,
    0.74458874458874458
)

$startCol = 12  # Column L

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $row2[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
}

$ws.Range("A1:S3").Select() | Out-Null
